# Populate Sheet3 with the bill computation data, then move the active
# sheet / selection from Sheet2 to Sheet3.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("G2").Value = 6870
$ws3.Range("G3").Value = 1616
$ws3.Range("A4").Value = 2082
$ws3.Range("G4").Formula = "=SUM(G2:G3)"
$ws3.Range("C5").Value = "RMB"
$ws3.Range("A9").Value = 1950
$ws3.Range("B9").Value = "2082-125"
$ws3.Range("C15").Value = 8486
$ws3.Range("D15").Value = "hkd"
$ws3.Range("A20").Value = 31514
$ws3.Range("B20").Value = "RNB"

# Sheet2 was the previously-active sheet/tab; move the selection before
# switching away so it keeps the new (non-active) selection state.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B37").Select()

# Make Sheet3 the active tab, with its own selection at D20.
$ws3.Select()
$ws3.Range("D20").Select()
